# Estadisticos Matutinos 15 Oct
# Fill in the first-partial (1P) and final grade statistics, propagate the
# "Reprobados/Aprobados" count into the 2nd-partial blanks column, and
# refresh the "Rescatables" (make-up exam) roster with the full list of
# students still pending a grade.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Estadisticos 1P": record Blancos/Reprobados/Aprobados/Por_Apro/
# Promedio for each of the 3 groups now that grades came in.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")

$ws1.Range("D2").Value = 0
$ws1.Range("E2").Value = 10
$ws1.Range("F2").Value = 21
$ws1.Range("G2").Value = 67.74
$ws1.Range("H2").Value = 7.1

$ws1.Range("D3").Value = 0
$ws1.Range("E3").Value = 9
$ws1.Range("F3").Value = 22
$ws1.Range("G3").Value = 70.97
$ws1.Range("H3").Value = 7.1

$ws1.Range("D4").Value = 0
$ws1.Range("E4").Value = 8
$ws1.Range("F4").Value = 26
$ws1.Range("G4").Value = 76.47
$ws1.Range("H4").Value = 7.3

# ---------------------------------------------------------------------
# Sheet "Estadisticos 2P": the "Reprobados" column now mirrors the group
# total (every student is still pending 2nd-partial grading).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")

$ws2.Range("E2").Value = 31
$ws2.Range("E3").Value = 31
$ws2.Range("E4").Value = 34

# ---------------------------------------------------------------------
# Sheet "Estadisticos Final": same figures as 1P (first-partial grades
# are currently the only ones available for the final computation).
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")

$ws3.Range("D2").Value = 0
$ws3.Range("E2").Value = 10
$ws3.Range("F2").Value = 21
$ws3.Range("G2").Value = 67.74
$ws3.Range("H2").Value = 7.1

$ws3.Range("D3").Value = 0
$ws3.Range("E3").Value = 9
$ws3.Range("F3").Value = 22
$ws3.Range("G3").Value = 70.97
$ws3.Range("H3").Value = 7.1

$ws3.Range("D4").Value = 0
$ws3.Range("E4").Value = 8
$ws3.Range("F4").Value = 26
$ws3.Range("G4").Value = 76.47
$ws3.Range("H4").Value = 7.3

# ---------------------------------------------------------------------
# Sheet "Rescatables": the roster grows from 2 students to 8 (6 new
# names inserted before the original 2, which move down to rows 8-9).
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")

# Insert 6 new rows above the existing data rows (old rows 2-3 become
# rows 8-9) before writing the new roster. Reset to the default "Normal"
# style so the inserted rows don't inherit the header's bold/border look.
$ws4.Range("A2:A7").EntireRow.Insert()
$ws4.Range("A2:G7").Style = "Normal"

$ws4.Range("A2").Value = 19330051920177
$ws4.Range("B2").Value = "SANCHEZ"
$ws4.Range("C2").Value = "MARTINEZ"
$ws4.Range("D2").Value = "JOSSELIN ANDREA"
$ws4.Range("E2").Value = "BIOLOGÍA"
$ws4.Range("F2").Value = "3ALCV"
$ws4.Range("G2").Value = 6

$ws4.Range("A3").Value = 19330051920177
$ws4.Range("B3").Value = "SANCHEZ"
$ws4.Range("C3").Value = "MARTINEZ"
$ws4.Range("D3").Value = "JOSSELIN ANDREA"
$ws4.Range("E3").Value = "IDENTIFICA MICROORGANISMOS CON BASE EN TÉCNICAS BACTERIOLÓGICAS"
$ws4.Range("F3").Value = "3ALCV"
$ws4.Range("G3").Value = 6

$ws4.Range("A4").Value = 19330051920273
$ws4.Range("B4").Value = "ARIAS"
$ws4.Range("C4").Value = "BARRAGAN"
$ws4.Range("D4").Value = "ESDRAS ALAN"
$ws4.Range("E4").Value = "ANALIZA Y FRACCIONA SANGRE CON FINES TRANSFUSIONALES"
$ws4.Range("F4").Value = "5ALCV"
$ws4.Range("G4").Value = 6

$ws4.Range("A5").Value = 19330051920278
$ws4.Range("B5").Value = "GARCIA"
$ws4.Range("C5").Value = "LINARES"
$ws4.Range("D5").Value = "ANDRES"
$ws4.Range("E5").Value = "ANALIZA Y FRACCIONA SANGRE CON FINES TRANSFUSIONALES"
$ws4.Range("F5").Value = "5ALCV"
$ws4.Range("G5").Value = 6

$ws4.Range("A6").Value = 19330051920286
$ws4.Range("B6").Value = "MANZANET"
$ws4.Range("C6").Value = "ANDRADE"
$ws4.Range("D6").Value = "JADE EMILY"
$ws4.Range("E6").Value = "ANALIZA Y FRACCIONA SANGRE CON FINES TRANSFUSIONALES"
$ws4.Range("F6").Value = "5ALCV"
$ws4.Range("G6").Value = 6

$ws4.Range("A7").Value = 18330051920346
$ws4.Range("B7").Value = "TEXCAHUA"
$ws4.Range("C7").Value = "CABRERA"
$ws4.Range("D7").Value = "YADIRA"
$ws4.Range("E7").Value = "ANALIZA Y FRACCIONA SANGRE CON FINES TRANSFUSIONALES"
$ws4.Range("F7").Value = "5ALCV"
$ws4.Range("G7").Value = 3

# Rows 8-9 keep the original two students, now with the updated
# "ANALIZA Y FRACCIONA..." subject/group/count for row 9.
$ws4.Range("E9").Value = "ANALIZA Y FRACCIONA SANGRE CON FINES TRANSFUSIONALES"
$ws4.Range("F9").Value = "5ALCV"
$ws4.Range("G9").Value = 2
